# Re-sort the comma-separated "Recorded By" names in column G
# (case-insensitive alphabetical order) for every data row on the
# active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $val = $cell.Value2

    if ($null -ne $val -and $val -ne "") {
        $parts = $val -split ",\s*"
        if ($parts.Count -gt 1) {
            $sorted = $parts | Sort-Object { $_.ToLower() }
            $newVal = [string]::Join(", ", $sorted)
            if ($newVal -ne $val) {
                $cell.Value2 = $newVal
            }
        }
    }
}
